# Auto-generated: update loading_percent values for the 380 kV case (Case_3_7)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 10.36309431898779
    "C2" = 6.008478398229731
    "D2" = 9.468811434381657
    "F2" = 36.79618150162513
    "G2" = 40.19051235771648
    "H2" = 17.15375422720692
    "I2" = 24.78883611322316
    "J2" = 11.24723541890542
    "K2" = 10.88095218673438
    "M2" = 16.79539424430797
    "N2" = 20.45882285062337
    "B3" = 10.1333009912424
    "C3" = 5.849672283108495
    "D3" = 9.436519386766914
    "F3" = 36.8240459191758
    "G3" = 40.20620052446105
    "H3" = 17.19580478590421
    "I3" = 24.85989192301721
    "J3" = 11.26420184048702
    "K3" = 10.73596941530452
    "M3" = 16.73433211890373
    "N3" = 20.52118954243674
    "B4" = 9.991849509082098
    "C4" = 5.751430126401396
    "D4" = 9.418246110746011
    "F4" = 36.84927329810952
    "G4" = 40.22654818640714
    "H4" = 17.22431316792005
    "I4" = 24.90790012352476
    "J4" = 11.27616933791612
    "K4" = 10.64826219304587
    "M4" = 16.69983547633253
    "N4" = 20.56128191324386
    "B5" = 9.934205135428362
    "C5" = 5.711276370208073
    "D5" = 9.411195820166467
    "F5" = 36.8615934638434
    "G5" = 40.23753067495275
    "H5" = 17.23660635902714
    "I5" = 24.92856372933726
    "J5" = 11.28143607634191
    "K5" = 10.61289392274797
    "M5" = 16.6865416526066
    "N5" = 20.57807350223735
    "B6" = 9.924635849884716
    "C6" = 5.704603635443942
    "D6" = 9.410049205520989
    "F6" = 36.86376235199712
    "G6" = 40.23951667699048
    "H6" = 17.23868843744754
    "I6" = 24.93206129248586
    "J6" = 11.28233416493823
    "K6" = 10.60704483095217
    "M6" = 16.68438064868356
    "N6" = 20.58088916746848
    "B7" = 9.991071985678472
    "C7" = 5.750888993119414
    "D7" = 9.418149416940174
    "F7" = 36.8494311959436
    "G7" = 40.22668541192056
    "H7" = 17.22447622252558
    "I7" = 24.90817434873058
    "J7" = 11.27623878823868
    "K7" = 10.64778363695697
    "M7" = 16.69965308513421
    "N7" = 20.56150653212963
    "B8" = 10.28398782191298
    "C8" = 5.95391280850754
    "D8" = 9.457357930783443
    "F8" = 36.80410339760385
    "G8" = 40.19369594606872
    "H8" = 17.16769488019919
    "I8" = 24.81242632035666
    "J8" = 11.25276380871212
    "K8" = 10.83071562063993
    "M8" = 16.7737248202589
    "N8" = 20.47995426987181
    "B9" = 10.85183346552037
    "C9" = 6.343503305270081
    "D9" = 9.546322739937999
    "F9" = 36.77967465069341
    "G9" = 40.21412277215409
    "H9" = 17.0777023122008
    "I9" = 24.6594775848618
    "J9" = 11.21902253958107
    "K9" = 11.1979525886179
    "M9" = 16.94225051864003
    "N9" = 20.33424610782242
    "B10" = 11.26029329683221
    "C10" = 6.621158859460929
    "D10" = 9.618682858654092
    "F10" = 36.80103281621984
    "G10" = 40.2810586296176
    "H10" = 17.02462742064742
    "I10" = 24.56841223356635
    "J10" = 11.2017179432068
    "K10" = 11.47039073079593
    "M10" = 17.07953522459797
    "N10" = 20.23577716226083
    "B11" = 11.44329656072012
    "C11" = 6.744975426653982
    "D11" = 9.653035754799468
    "F11" = 36.81926750923235
    "G11" = 40.32275742305897
    "H11" = 17.00331879527393
    "I11" = 24.53162745002153
    "J11" = 11.1954681491036
    "K11" = 11.5943658374798
    "M11" = 17.1447461735666
    "N11" = 20.19282700120401
    "B12" = 11.51212124973784
    "C12" = 6.791456213662187
    "D12" = 9.666243078509524
    "F12" = 36.82739435718872
    "G12" = 40.34015986688264
    "H12" = 16.99565780565137
    "I12" = 24.51836672323146
    "J12" = 11.19333444373267
    "K12" = 11.64127438001527
    "M12" = 17.16982143577661
    "N12" = 20.17682673570284
    "B13" = 11.49732081112666
    "C13" = 6.781464511767808
    "D13" = 9.663389939196456
    "F13" = 36.82558981851449
    "G13" = 40.33634035331994
    "H13" = 16.99728957959326
    "I13" = 24.52119288729061
    "J13" = 11.19378362005131
    "K13" = 11.63117417224584
    "M13" = 17.16440432463455
    "N13" = 20.18026095406297
    "B14" = 11.44896871174677
    "C14" = 6.748807821869129
    "D14" = 9.654118394508531
    "F14" = 36.81991164586554
    "G14" = 40.32415686407864
    "H14" = 17.00268033952754
    "I14" = 24.53052306875794
    "J14" = 11.19528794104611
    "K14" = 11.59822604441166
    "M14" = 17.14680157978933
    "N14" = 20.19150536455897
    "B15" = 11.41928781449196
    "C15" = 6.728750454664759
    "D15" = 9.648464934710708
    "F15" = 36.81659259864303
    "G15" = 40.31690385914939
    "H15" = 17.00603549642639
    "I15" = 24.53632522326827
    "J15" = 11.19623970879495
    "K15" = 11.57803810833974
    "M15" = 17.13606857514743
    "N15" = 20.19842724203881
    "B16" = 11.2482709338505
    "C16" = 6.613012867223736
    "D16" = 9.616466045666289
    "F16" = 36.80001230245721
    "G16" = 40.27855944152747
    "H16" = 17.02607708002463
    "I16" = 24.57090972692623
    "J16" = 11.20215899874625
    "K16" = 11.46228570549518
    "M16" = 17.0753277503599
    "N16" = 20.23862106371916
    "B17" = 11.14258655902845
    "C17" = 6.541338480311881
    "D17" = 9.597198042237071
    "F17" = 36.79202070671429
    "G17" = 40.25791418413615
    "H17" = 17.03909845523774
    "I17" = 24.59331600970674
    "J17" = 11.2062055618733
    "K17" = 11.39125426629552
    "M17" = 17.03876113681252
    "N17" = 20.26375012933637
    "B18" = 11.08153941709587
    "C18" = 6.499881456899582
    "D18" = 9.586251241247446
    "F18" = 36.78822623298178
    "G18" = 40.24709887106479
    "H18" = 17.0468548671337
    "I18" = 24.60664028942232
    "J18" = 11.20868574448965
    "K18" = 11.35040535441375
    "M18" = 17.01798968851728
    "N18" = 20.27837732158396
    "B19" = 11.06082744353631
    "C19" = 6.485806530047481
    "D19" = 9.582568374130579
    "F19" = 36.78707932956557
    "G19" = 40.24361906809989
    "H19" = 17.04952687580256
    "I19" = 24.61122663309418
    "J19" = 11.20955172770534
    "K19" = 11.3365770930451
    "M19" = 17.01100207144879
    "N19" = 20.28335969089261
    "B20" = 11.1538643095915
    "C20" = 6.54899269182765
    "D20" = 9.599235171318506
    "F20" = 36.79278843354457
    "G20" = 40.26000230910795
    "H20" = 17.0376846854418
    "I20" = 24.59088560689382
    "J20" = 11.20575899582556
    "K20" = 11.39881533553504
    "M20" = 17.04262684418948
    "N20" = 20.26105713666153
    "B21" = 11.46318431946603
    "C21" = 6.758411258333963
    "D21" = 9.656836345283494
    "F21" = 36.82154633501359
    "G21" = 40.32769175440596
    "H21" = 17.00108586336051
    "I21" = 24.52776440403087
    "J21" = 11.19483976586369
    "K21" = 11.60790508910832
    "M21" = 17.15196171106832
    "N21" = 20.18819545180543
    "B22" = 11.66254660053295
    "C22" = 6.892893476721837
    "D22" = 9.695635951946187
    "F22" = 36.84746047859357
    "G22" = 40.38132288883387
    "H22" = 16.97954542847979
    "I22" = 24.49041037275475
    "J22" = 11.18906108296702
    "K22" = 11.7443154354895
    "M22" = 17.22563322532379
    "N22" = 20.1421145790063
    "B23" = 11.55642123950658
    "C23" = 6.821350720701449
    "D23" = 9.674824957884683
    "F23" = 36.83297952525656
    "G23" = 40.35184188720486
    "H23" = 16.99082416577406
    "I23" = 24.50998967134343
    "J23" = 11.19202116390429
    "K23" = 11.67154706119211
    "M23" = 17.18611596616184
    "N23" = 20.16656840527764
    "B24" = 11.14876652949144
    "C24" = 6.545532999263179
    "D24" = 9.598313777775406
    "F24" = 36.792438852024
    "G24" = 40.25905498446061
    "H24" = 17.03832300895935
    "I24" = 24.59198301338209
    "J24" = 11.20596040936494
    "K24" = 11.39539701009235
    "M24" = 17.04087837429245
    "N24" = 20.26227407735676
    "B25" = 10.6994233838182
    "C25" = 6.239386108474565
    "D25" = 9.520998669071187
    "F25" = 36.77937778908152
    "G25" = 40.19947345237782
    "H25" = 17.09975882681433
    "I25" = 24.69711848549695
    "J25" = 11.2268351344234
    "K25" = 11.09795956841579
    "M25" = 16.8942396147789
    "N25" = 20.37215081225482
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

